# Scheduled-runner market data refresh.
# Updates currentAveragePrice / currentAveragePriceNQ/HQ / LevePriceNQ/HQ /
# LeveProfitNQ/HQ (columns H:N) for the affected leve rows on each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 84.666664
$ws.Range("I6").Value = 84.666664
$ws.Range("K6").Value = 253.999992
$ws.Range("M6").Value = -141.999992

$ws.Range("H28").Value = 2334.8333
$ws.Range("I28").Value = 811.2
$ws.Range("K28").Value = 811.2
$ws.Range("M28").Value = -326.2

$ws.Range("H53").Value = 775.5714
$ws.Range("J53").Value = 250
$ws.Range("L53").Value = 250
$ws.Range("N53").Value = -1524

$ws.Range("H76").Value = 9451.666999999999
$ws.Range("I76").Value = 6750
$ws.Range("J76").Value = 9992
$ws.Range("K76").Value = 6750
$ws.Range("L76").Value = 9992
$ws.Range("M76").Value = -6435
$ws.Range("N76").Value = -10622

$ws.Range("H79").Value = 9451.666999999999
$ws.Range("I79").Value = 6750
$ws.Range("J79").Value = 9992
$ws.Range("K79").Value = 6750
$ws.Range("L79").Value = 9992
$ws.Range("M79").Value = -5658
$ws.Range("N79").Value = -12176

$ws.Range("H98").Value = 4316.5835
$ws.Range("I98").Value = 4316.5835
$ws.Range("K98").Value = 4316.5835
$ws.Range("M98").Value = -2818.5835

$ws.Range("H106").Value = 4076
$ws.Range("I106").Value = 4076
$ws.Range("K106").Value = 4076
$ws.Range("M106").Value = -3445

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H112").Value = 2435.7896
$ws.Range("J112").Value = 2435.7896
$ws.Range("L112").Value = 7307.3688
$ws.Range("N112").Value = -9523.3688

$ws.Range("H122").Value = 4316.5835
$ws.Range("I122").Value = 4316.5835
$ws.Range("K122").Value = 12949.7505
$ws.Range("M122").Value = -10499.7505

$ws.Range("H126").Value = 77740
$ws.Range("J126").Value = 77740
$ws.Range("L126").Value = 77740
$ws.Range("N126").Value = -87620

$ws.Range("H129").Value = 1721.5862
$ws.Range("J129").Value = 3875.7
$ws.Range("L129").Value = 11627.1
$ws.Range("N129").Value = -21627.1

$ws.Range("H132").Value = 6169.8887
$ws.Range("I132").Value = 5035.8237
$ws.Range("K132").Value = 15107.4711
$ws.Range("M132").Value = -12577.4711

$ws.Range("H138").Value = 2980.0256
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2980.0256
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 8940.076799999999
$ws.Range("N138").Value = -19220.0768
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2562.35
$ws.Range("I32").Value = 2562.35
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2562.35
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2275.35
$ws.Range("N32").ClearContents()

$ws.Range("H122").Value = 4231.3667
$ws.Range("I122").Value = 4149.115
$ws.Range("K122").Value = 12447.345
$ws.Range("M122").Value = -9997.344999999999

$ws.Range("H124").Value = 60713.5
$ws.Range("J124").Value = 60713.5
$ws.Range("L124").Value = 60713.5
$ws.Range("N124").Value = -70533.5

$ws.Range("H132").Value = 2537.5925
$ws.Range("I132").Value = 1405.5238
$ws.Range("K132").Value = 4216.5714
$ws.Range("M132").Value = -1686.5714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 181819220
$ws.Range("J94").Value = 618.5
$ws.Range("L94").Value = 618.5
$ws.Range("N94").Value = -1520.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4035.149
$ws.Range("I31").Value = 2920.1853
$ws.Range("J31").Value = 5540.35
$ws.Range("K31").Value = 2920.1853
$ws.Range("L31").Value = 5540.35
$ws.Range("M31").Value = -2625.1853
$ws.Range("N31").Value = -6130.35

$ws.Range("H34").Value = 4035.149
$ws.Range("I34").Value = 2920.1853
$ws.Range("J34").Value = 5540.35
$ws.Range("K34").Value = 2920.1853
$ws.Range("L34").Value = 5540.35
$ws.Range("M34").Value = -2718.1853
$ws.Range("N34").Value = -5944.35

$ws.Range("H58").Value = 2866.5518
$ws.Range("I58").Value = 2413.077
$ws.Range("J58").Value = 3235
$ws.Range("K58").Value = 2413.077
$ws.Range("L58").Value = 3235
$ws.Range("M58").Value = -2210.077
$ws.Range("N58").Value = -3641

$ws.Range("H106").Value = 67018.39999999999
$ws.Range("J106").Value = 67018.39999999999
$ws.Range("L106").Value = 67018.39999999999
$ws.Range("N106").Value = -69542.39999999999

$ws.Range("H122").Value = 4398.615
$ws.Range("J122").Value = 5705.8335
$ws.Range("L122").Value = 17117.5005
$ws.Range("N122").Value = -22017.5005

$ws.Range("H132").Value = 3687.4583
$ws.Range("I132").Value = 3728.7
$ws.Range("K132").Value = 11186.1
$ws.Range("M132").Value = -8656.099999999999

$ws.Range("H136").Value = 2866.5518
$ws.Range("I136").Value = 2413.077
$ws.Range("J136").Value = 3235
$ws.Range("K136").Value = 7239.231000000001
$ws.Range("L136").Value = 9705
$ws.Range("M136").Value = -4689.231000000001
$ws.Range("N136").Value = -14805

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 996
$ws.Range("I8").Value = 996
$ws.Range("K8").Value = 2988
$ws.Range("M8").Value = -2849

$ws.Range("H121").Value = 500450
$ws.Range("I121").Value = 1000150
$ws.Range("J121").Value = 250600
$ws.Range("K121").Value = 3000450
$ws.Range("L121").Value = 751800
$ws.Range("M121").Value = -2999140
$ws.Range("N121").Value = -754420

$ws.Range("H132").Value = 2103.3157
$ws.Range("J132").Value = 3485.4285
$ws.Range("L132").Value = 31368.8565
$ws.Range("N132").Value = -36428.8565

$ws.Range("H134").Value = 1426.0769
$ws.Range("I134").Value = 1426.0769
$ws.Range("K134").Value = 4278.2307
$ws.Range("M134").Value = 791.7692999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 668666.7
$ws.Range("I70").Value = 668666.7
$ws.Range("K70").Value = 668666.7
$ws.Range("M70").Value = -668396.7

$ws.Range("H73").Value = 668666.7
$ws.Range("I73").Value = 668666.7
$ws.Range("K73").Value = 668666.7
$ws.Range("M73").Value = -667730.7

$ws.Range("H80").Value = 90912770
$ws.Range("J80").Value = 3927.2856
$ws.Range("L80").Value = 3927.2856
$ws.Range("N80").Value = -5923.2856

$ws.Range("H83").Value = 90912770
$ws.Range("J83").Value = 3927.2856
$ws.Range("L83").Value = 19636.428
$ws.Range("N83").Value = -29620.428

$ws.Range("H113").Value = 5333
$ws.Range("I113").Value = 4999.5
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 4999.5
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -2829.5
$ws.Range("N113").Value = -10340

$ws.Range("H122").Value = 6999.1665
$ws.Range("I122").Value = 5248.75
$ws.Range("K122").Value = 15746.25
$ws.Range("M122").Value = -13296.25

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H126").Value = 19702.666
$ws.Range("I126").Value = 17554.5
$ws.Range("J126").Value = 23999
$ws.Range("K126").Value = 52663.5
$ws.Range("L126").Value = 71997
$ws.Range("M126").Value = -50193.5
$ws.Range("N126").Value = -76937

$ws.Range("H133").Value = 122483.75
$ws.Range("J133").Value = 122483.75
$ws.Range("L133").Value = 122483.75
$ws.Range("N133").Value = -132603.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 947.25
$ws.Range("I16").Value = 939.7143
$ws.Range("K16").Value = 939.7143
$ws.Range("M16").Value = -769.7143

$ws.Range("H22").Value = 1029.5
$ws.Range("J22").Value = 1164.3334
$ws.Range("L22").Value = 1164.3334
$ws.Range("N22").Value = -1754.3334

$ws.Range("H27").Value = 1029.5
$ws.Range("J27").Value = 1164.3334
$ws.Range("L27").Value = 1164.3334
$ws.Range("N27").Value = -1378.3334

$ws.Range("H61").Value = 2549.55
$ws.Range("J61").Value = 2728
$ws.Range("L61").Value = 2728
$ws.Range("N61").Value = -3132

$ws.Range("H68").Value = 3998
$ws.Range("I68").Value = 3998
$ws.Range("K68").Value = 3998
$ws.Range("M68").Value = -3249

$ws.Range("H71").Value = 3998
$ws.Range("I71").Value = 3998
$ws.Range("K71").Value = 19990
$ws.Range("M71").Value = -16246

$ws.Range("H100").Value = 3876.7778
$ws.Range("I100").Value = 4256
$ws.Range("J100").Value = 2549.5
$ws.Range("K100").Value = 4256
$ws.Range("L100").Value = 2549.5
$ws.Range("M100").Value = -3715
$ws.Range("N100").Value = -3631.5

$ws.Range("H104").Value = 25000
$ws.Range("J104").Value = 25000
$ws.Range("L104").Value = 25000
$ws.Range("N104").Value = -31988

$ws.Range("H113").Value = 2549.55
$ws.Range("J113").Value = 2728
$ws.Range("L113").Value = 2728
$ws.Range("N113").Value = -7068

$ws.Range("H132").Value = 4746.6553
$ws.Range("I132").Value = 2841.8333
$ws.Range("J132").Value = 7863.636
$ws.Range("K132").Value = 8525.499899999999
$ws.Range("L132").Value = 23590.908
$ws.Range("M132").Value = -5995.499899999999
$ws.Range("N132").Value = -28650.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H104").Value = 20545.9
$ws.Range("J104").Value = 20545.9
$ws.Range("L104").Value = 20545.9
$ws.Range("N104").Value = -27533.9

$ws.Range("H132").Value = 7833.3335
$ws.Range("I132").Value = 9900
$ws.Range("K132").Value = 29700
$ws.Range("M132").Value = -27170
